# Mark row 24 ("25/08/2023" / "Add User form" / "grid data is not shown on
# controls esp textboxes") as completed: copy the "done" formatting
# (strikethrough font, etc.) from an already-completed row onto row 24, and
# stamp the completion date into column D.
#
# The row that used to be blank (row 25) becomes the new pending task:
# "25/08/2023" / "Add User form" / "password displayed on the control should
# not show the characters".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the "task completed" formatting (font/fill/border/alignment) from an
# existing completed row (row 20) onto row 24 so it gets struck through like
# the other finished tasks.
$ws.Range("A20:D20").Copy() | Out-Null
$ws.Range("A24:D24").PasteSpecial(-4122) | Out-Null

# Record when the task on row 24 was completed.
$ws.Range("D24").Value = "25/08/23"

# Add the new task entry on row 25 (formatting already matches the other
# pending rows, so only the values need to be filled in).
$ws.Range("A25").Value = "25/08/2023"
$ws.Range("B25").Value = "Add User form"
$ws.Range("C25").Value = "password displayed on the control should not show the characters"

# Reflect the author's new selection.
$ws.Range("C25").Select() | Out-Null
